# IMU and XBee blocks assigned to XCore 3.
# - Signals set.
# - Initial off-board trace routing.

$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item(5)   # "XMOS Dualchip"
$ws6 = $wb.Worksheets.Item(6)   # "XMOS Dualchip Planning"

# ---------------------------------------------------------------------------
# Sheet "XMOS Dualchip Planning" (ws6): new "Unallocated" column (L) on the
# by-core table, a couple of individual cell fill-ins, the new totals row,
# and a new "1-bit ports remaining" summary row.
# ---------------------------------------------------------------------------

# New unique strings are appended to the shared-string table in the order
# cells are first written, so write the brand-new strings in source order:
# '1-bit ports remaining', 'Unallocated', 'U13_DRDY', 'XBEE_CTS', 'XBEE_RTS',
# 'U15_INT1', 'U15_INT2'.
$ws6.Range("G14").Value = "1-bit ports remaining"
$ws6.Range("L1").Value = "Unallocated"

# Per-block "Unallocated" 1-bit port counts for Pmod1, Pmod2, Gadgeteer
# (new column L), plus the per-core allocations for I2C/STM32 UART/Xbee UART
# -- write these BEFORE the totals formulas below so the totals recalc off
# live numbers instead of stale cached ones.
$ws6.Range("L3").Value = 8
$ws6.Range("L4").Value = 8
$ws6.Range("L5").Value = 8

# I2C block: 3 remaining 1-bit ports went to core X2 (column K).
$ws6.Range("K7").Value = 3

# STM32 UART / Xbee UART remainders.
$ws6.Range("L10").Value = 3
$ws6.Range("K11").Value = 4

# Totals row 13: the shared formula spanning I13:K13 effectively grows to
# I13:L13 (K13 = SUM(K2:K11) = 16, new L13 = SUM(L2:L11) = 27).
$ws6.Range("K13").Formula = "=SUM(K2:K11)"
$ws6.Range("L13").Formula = "=SUM(L2:L11)"

# New row 14: 1-bit ports remaining per core (16 total 1-bit ports/core).
$ws6.Range("H14").Formula = "=16-H13"
$ws6.Range("I14").Formula = "=16-I13"
$ws6.Range("J14").Formula = "=16-J13"
$ws6.Range("K14").Formula = "=16-K13"

# ---------------------------------------------------------------------------
# Sheet "XMOS Dualchip" (ws5): signal assignments for the newly-placed IMU
# (U13/U14/U15) and XBee blocks on XCore 3 (columns N:P).
# ---------------------------------------------------------------------------

# U13_DRDY / XBEE_CTS / XBEE_RTS are brand-new shared strings, so write them
# ahead of the remaining (already-reused) ones to keep shared-string order
# matching the source.
$ws5.Range("P27").Value = "U13_DRDY"
$ws5.Range("P40").Value = "XBEE_CTS"
$ws5.Range("P41").Value = "XBEE_RTS"

$ws5.Range("P2").Value = "SDA"
$ws5.Range("P3").Value = "SCL"

$ws5.Range("O4").Value = "P4A0"
$ws5.Range("P4").Value = "U15_INT1"

$ws5.Range("O5").Value = "P4A1"
$ws5.Range("P5").Value = "U15_INT2"

$ws5.Range("O6").Value = "P4B0"
$ws5.Range("P6").Value = "U14_INT2"

$ws5.Range("O7").Value = "P4B1"
$ws5.Range("P7").Value = "U14_INT1"

$ws5.Range("O8").Value = "P4B2"
$ws5.Range("O9").Value = "P4B3"
$ws5.Range("O10").Value = "P4A2"
$ws5.Range("O11").Value = "P4A3"

$ws5.Range("P12").Value = "XBEE_DOUT"
$ws5.Range("P13").Value = "XBEE_DIN"

# ---------------------------------------------------------------------------
# View state: XMOS Dualchip becomes the active/selected sheet & tab, scrolled
# so column E is leftmost with P8 selected; XMOS Dualchip Planning drops its
# tab selection and its own selection moves to J17.
# ---------------------------------------------------------------------------

$ws6.Range("J17").Select() | Out-Null
$ws5.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws5.Range("P8").Select() | Out-Null
